# "correção de erros, implementação de data de criação"
# Adds a new column Z ("Data de criação") that is derived from the date
# suffix present at the end of each ad's name (column A), e.g.
#   "AD 01 - CARROSSEL FEED - 05/06/2024"  ->  "05/06/2024"
# Rows whose name has no trailing date (or a 2-digit year, e.g. "04/04/24")
# are left blank in column Z.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 33

# --- New header cell Z1 ("Data de criação"), matching the style of the
#     existing X1/Y1 report-range headers -----------------------------------
$ws.Range("X1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Z1").Value = "Data de criação"

# --- Force columns X:Z (rows 2..lastRow) to Text format so that the
#     dd/mm/yyyy strings are never silently reinterpreted as date serials ---
$ws.Range("X2:Z" + $lastRow).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    # Correção: garante que o período do relatório está preenchido em todas as linhas
    $ws.Cells.Item($r, 24).Value = "01/01/2024"   # X: Início dos relatórios
    $ws.Cells.Item($r, 25).Value = "30/10/2024"   # Y: Término dos relatórios

    # Implementação: extrai a data de criação a partir do nome do anúncio
    $adName = $ws.Cells.Item($r, 1).Value2
    if ($adName -match '(\d{2}/\d{2}/\d{4})\s*$') {
        $ws.Cells.Item($r, 26).Value = $matches[1]   # Z: Data de criação
    } else {
        $ws.Cells.Item($r, 26).Value = ""
    }
}
